{"js": "// Applies the \u0422\u0417.docx revision:\n//  - tightens the W / H measurement ranges\n//  - replaces the \"D\" bullet's free-text upper bound with a formula (D = W/3 \u043c\u043c)\n//    and drops the reviewer comment attached to it\n//  - turns the \"*\" multiplications in the T bullet into \"/\" divisions and\n//    drops the reviewer comment attached to it\n//  - rewrites the N formula and drops the reviewer comment attached to it\n//  - fixes the N bullet's left indent to match its first-line indent\n\n// 1) Delete every comment in the document (also strips the now-orphaned\n//    commentRangeStart/End + commentReference runs from the body).\nconst comments = context.document.body.getComments();\ncomments.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < comments.items.length; i++) {\n  comments.items[i].delete();\n}\nawait context.sync();\n\n// 2) Small, targeted text substitutions (keeps surrounding run formatting\n//    intact as much as the search/replace API allows).\nconst replacements = [\n  [\"10 \u2013 100\", \"15 \u2013 100\"],\n  [\"\u0440\u0430\u0437\u043c\u0435\u0440\u043e\u0432 \u043e\u0433\u0440\u0430\u043d\u0438\u0447\u0435\u043d\u043d\u044b\u0445 \u043c\u043e\u0434\u0435\u043b\u044c\u044e);\", \"D = W/3 \u043c\u043c); \"],\n  [\"50 \u2013 200\", \"15 \u2013 150\"],\n  [\"H*6, T<W*3.)\", \"H/6, T<W/6)\"],\n  [\"N=(H-L-M)/(D+5));\", \"N = W / (D * 1,5));\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Note: the \"D\" bullet's closing \");\" needs to collapse to \"); \" (the\n// surviving \")\" / \";\" runs stay untouched, but the paragraph needs to pick\n// up the trailing space that followed the deleted comment marker in the\n// revision). The substring replace above already appended it.\n\n// 3) Fix the left indent of the \"N \u2013 ...\" bullet (was 0, should match the\n//    851-twip first-line indent, i.e. 42.55 points).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"\u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e\") !== -1 && p.text.indexOf(\"N\") === 0) {\n    p.leftIndent = 42.55;\n  }\n}\nawait context.sync();\n", "ps1": "# Applies the \u0422\u0417.docx revision:\n#  - tightens the W / H measurement ranges\n#  - replaces the \"D\" bullet's free-text upper bound with a formula (D = W/3 \u043c\u043c)\n#    and drops the reviewer comment attached to it\n#  - turns the \"*\" multiplications in the T bullet into \"/\" divisions and\n#    drops the reviewer comment attached to it\n#  - rewrites the N formula and drops the reviewer comment attached to it\n#  - fixes the N bullet's left indent to match its first-line indent\n\n$d = $word.ActiveDocument\n\n# 1) Delete every comment in the document (also strips the now-orphaned\n#    commentRangeStart/End + commentReference runs from the body).\nfor ($i = $d.Comments.Count; $i -ge 1; $i--) {\n    $d.Comments.Item($i).Delete()\n}\n\n# 2) Small, targeted text substitutions (keeps surrounding run formatting\n#    intact as much as the Find/Replace engine allows).\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text \"W \u2013 \u0448\u0438\u0440\u0438\u043d\u0430 \u0443\u0433\u043e\u043b\u043a\u0430 (10 \u2013 100 \u043c\u043c)\" \"W \u2013 \u0448\u0438\u0440\u0438\u043d\u0430 \u0443\u0433\u043e\u043b\u043a\u0430 (15 \u2013 100 \u043c\u043c)\"\nReplace-Text \"\u0438 \u0434\u043e \u0440\u0430\u0437\u043c\u0435\u0440\u043e\u0432 \u043e\u0433\u0440\u0430\u043d\u0438\u0447\u0435\u043d\u043d\u044b\u0445 \u043c\u043e\u0434\u0435\u043b\u044c\u044e);\" \"\u0438 \u0434\u043e D = W/3 \u043c\u043c); \"\nReplace-Text \"H \u2013 \u0432\u044b\u0441\u043e\u0442\u0430 \u0443\u0433\u043e\u043b\u043a\u0430 (50 \u2013 200 \u043c\u043c)\" \"H \u2013 \u0432\u044b\u0441\u043e\u0442\u0430 \u0443\u0433\u043e\u043b\u043a\u0430 (15 \u2013 150 \u043c\u043c)\"\nReplace-Text \"T<H*6, T<W*3.)\" \"T<H/6, T<W/6)\"\nReplace-Text \"N=(H-L-M)/(D+5));\" \"N = W / (D * 1,5));\"\n\n# 3) Fix the left indent of the \"N \u2013 ...\" bullet (was 0, should match the\n#    851-twip first-line indent, i.e. 42.55 points).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"N\") -and $t.Contains(\"\u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e\")) {\n        $p.Format.LeftIndent = 42.55\n    }\n}\n"}
